$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9991625547409058
$ws.Range("B1").Value = 1.482583284378052
$ws.Range("D1").Value = 1.743406891822815
$ws.Range("E1").Value = 1.043423652648926
